# Apply "想去人数" (F) and "最低票价" (G) value updates across the four
# sheets of the workbook, matching the published gh-pages data refresh
# (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G5").Value  = 60
$ws1.Range("F7").Value  = 205
$ws1.Range("F10").Value = 432
$ws1.Range("F12").Value = 273
$ws1.Range("F13").Value = 99
$ws1.Range("F17").Value = 6334
$ws1.Range("F19").Value = 62
$ws1.Range("F21").Value = 7303
$ws1.Range("F24").Value = 3317
$ws1.Range("F25").Value = 442
$ws1.Range("F26").Value = 821
$ws1.Range("F27").Value = 4484
$ws1.Range("F28").Value = 338
$ws1.Range("F29").Value = 168
$ws1.Range("F30").Value = 165
$ws1.Range("F31").Value = 1338
$ws1.Range("F32").Value = 129
$ws1.Range("F35").Value = 1052
$ws1.Range("F36").Value = 1416
$ws1.Range("F37").Value = 2093

# --- Sheet 2: 演出 --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 57

# --- Sheet 3: 本地生活 -----------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 233
$ws3.Range("F3").Value = 1178
$ws3.Range("F4").Value = 62

# --- Sheet 4: 全部类型 -----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 233
$ws4.Range("F4").Value  = 1178
$ws4.Range("F5").Value  = 62
$ws4.Range("G8").Value  = 60
$ws4.Range("F10").Value = 205
$ws4.Range("F13").Value = 432
$ws4.Range("F15").Value = 273
$ws4.Range("F16").Value = 57
$ws4.Range("F17").Value = 99
$ws4.Range("F21").Value = 6334
$ws4.Range("F23").Value = 62
$ws4.Range("F25").Value = 7303
$ws4.Range("F28").Value = 3317
$ws4.Range("F29").Value = 442
$ws4.Range("F30").Value = 821
$ws4.Range("F31").Value = 4484
$ws4.Range("F32").Value = 338
$ws4.Range("F34").Value = 168
$ws4.Range("F35").Value = 165
$ws4.Range("F36").Value = 1338
$ws4.Range("F37").Value = 129
$ws4.Range("F40").Value = 1052
$ws4.Range("F41").Value = 1416
$ws4.Range("F43").Value = 2093
